$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header cells to their non-accented versions
$ws.Range("E1").Value = "Numero"
$ws.Range("F1").Value = "Dias"
$ws.Range("H1").Value = "Clinica"

# Update the active selection to H1, matching the saved sheet view state
$ws.Range("H1").Select()
